# PutShipCommand sequence diagram update
# ---------------------------------------------------------------
# Helper: find a shape on a Shapes collection by its PowerPoint shape Id.
function Get-ShapeById($shapes, $targetId) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $cand = $shapes.Item($i)
        if ($cand.Id -eq $targetId) { return $cand }
    }
    return $null
}

# Helper: replace the first occurrence of $old with $new inside a shape's
# text, touching only those characters so surrounding runs/formatting stay
# intact.
function Replace-InShapeText($shape, $old, $new) {
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($old)
    $rng = $tr.Characters($idx + 1, $old.Length)
    $rng.Text = $new
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# -----------------------------------------------------------------
# 1) "args" -> "args1" (two call-out labels on the left of the diagram)
# -----------------------------------------------------------------
$shParseCommand = Get-ShapeById $s.Shapes 80   # TextBox 79: parseCommand("put <args>")
Replace-InShapeText $shParseCommand "args" "args1"

$shParse = Get-ShapeById $s.Shapes 72          # TextBox 71: parse("<args>")
Replace-InShapeText $shParse "args" "args1"

# -----------------------------------------------------------------
# 2) "playerMapGrid" -> "mapGrid"
# -----------------------------------------------------------------
$shMapGrid = Get-ShapeById $s.Shapes 114       # TextBox 113
Replace-InShapeText $shMapGrid "playerMapGrid" "mapGrid"

# -----------------------------------------------------------------
# 3) "performChecks()" label repositioned
# -----------------------------------------------------------------
$shPerformChecks = Get-ShapeById $s.Shapes 129 # TextBox 128
$shPerformChecks.Left = 14053335 / 12700.0
$shPerformChecks.Top = 6010269 / 12700.0

# -----------------------------------------------------------------
# 4) Ungroup "Group 105" (id 106) -- its 4 children (136-139) become
#    top level shapes in the same spot in the z-order.
# -----------------------------------------------------------------
$grp = Get-ShapeById $s.Shapes 106
$grp.Ungroup() | Out-Null

# 4a) Curved Connector 135 (136) and Rectangle 136 (137) and
#     Curved Connector 12 (138) keep the position Ungroup already computed.

# 4b) TextBox 138 (139): "deployBattleship()" -> "deployBattleship(bs, cs, ot)"
#     and moved/resized to make room for the longer label.
$shDeployBattleship = Get-ShapeById $s.Shapes 139
Replace-InShapeText $shDeployBattleship "()" "(bs, cs, ot)"
$shDeployBattleship.Left = 10810081 / 12700.0
$shDeployBattleship.Top = 7077438 / 12700.0
$shDeployBattleship.Width = 1157705 / 12700.0
$shDeployBattleship.Height = 184666 / 12700.0

# -----------------------------------------------------------------
# 5) "isEnoughBattleships(battleship, 1)" -> "isEnoughBattleships(bs, 1)"
#    and its label repositioned.
# -----------------------------------------------------------------
$shIsEnough = Get-ShapeById $s.Shapes 105      # TextBox 104
Replace-InShapeText $shIsEnough "battleship" "bs"
$shIsEnough.Left = 8309745 / 12700.0
$shIsEnough.Top = 4953556 / 12700.0

# -----------------------------------------------------------------
# 6) Three new call-out text boxes at the end of the shape tree.
# -----------------------------------------------------------------

# 6a) id 93 "TextBox 92": " getFleet().deployOneBattleship(bs, cs, ot); "
#     Duplicated from the just-edited deployBattleship box (139) because it
#     already carries the matching run/endParaRPr colour scheme.
$dup93 = $shDeployBattleship.Duplicate().Item(1)
$tr93 = $dup93.TextFrame.TextRange
$full93 = $tr93.Text
$whole93 = $tr93.Characters(1, $full93.Length)
$whole93.Text = " getFleet().deployOneBattleship(bs, cs, ot); "
$dup93.Left = 10231396 / 12700.0
$dup93.Top = 7082490 / 12700.0
$dup93.Width = 3281404 / 12700.0
$dup93.Height = 184666 / 12700.0

# 6b) id 98 "TextBox 97": "bs = battleship, cs = coordinates, ot = orientation"
#     Duplicated from the args1/args2 legend box (61), which already has the
#     matching no-lstStyle / no-line styling.
$srcLegend = Get-ShapeById $s.Shapes 61
$dup98 = $srcLegend.Duplicate().Item(1)
$tr98 = $dup98.TextFrame.TextRange
$full98 = $tr98.Text
$whole98 = $tr98.Characters(1, $full98.Length)
$whole98.Text = "bs = battleship, cs = coordinates, ot = orientation"
$dup98.Left = 9864807 / 12700.0
$dup98.Top = 9172940 / 12700.0
$dup98.Width = 9254858 / 12700.0
$dup98.Height = 215444 / 12700.0

# 6c) id 101 "TextBox 100": "new BoundaryValueChecker(mapGrid, bs, cs, ot)"
#     Duplicated from the putShip() box (68), which already carries the
#     matching left-aligned paragraph with a coloured run.
$srcPutShip = Get-ShapeById $s.Shapes 68
$dup101 = $srcPutShip.Duplicate().Item(1)
$tr101 = $dup101.TextFrame.TextRange
$full101 = $tr101.Text
$whole101 = $tr101.Characters(1, $full101.Length)
$whole101.Text = "new BoundaryValueChecker(mapGrid, bs, cs, ot)"
$dup101.TextFrame.TextRange.Font.Color.RGB = RGB(0xA0, 0x30, 0x70)
$dup101.Left = 8287397 / 12700.0
$dup101.Top = 5461793 / 12700.0
$dup101.Width = 4440061 / 12700.0
$dup101.Height = 184666 / 12700.0
